$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Adiciona coluna "Ativo" com a logica de cadastro (S/N)
$ws.Range("C1").Value = "Ativo"
$ws.Range("C2").Value = "S"
$ws.Range("C3").Value = "S"
$ws.Range("C4").Value = "N"

$ws.Range("C5").Select()
